$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sectorData = @(
    @("Mortgage Real Estate Investment Trust...(16)", 0.5458030839494004),
    @("Banks(251)", 0.5109315629143277),
    @("Airlines(14)", 0.4876383660172545),
    @("Auto Components(21)", 0.4735946512814284),
    @("Leisure Products(11)", 0.4685931722379817),
    @("Multi-Utilities(18)", 0.4491987284344434),
    @("Trading Companies & Distributors(25)", 0.4447503474417077),
    @("Machinery(86)", 0.4190249390402377),
    @("Marine(15)", 0.4040424720855461),
    @("Building Products(24)", 0.3934823903959807),
    @("Energy Equipment & Services(38)", 0.3930019441961127),
    @("Road & Rail(22)", 0.372095315183784),
    @("Gas Utilities(12)", 0.3537754031430047),
    @("Semiconductors & Semiconductor Equipment(70)", 0.347408996185982),
    @("Containers & Packaging(12)", 0.342474928558079),
    @("Electric Utilities(28)", 0.335478499178457),
    @("Chemicals(52)", 0.3346378304215818),
    @("Textiles, Apparel & Luxury Goods(29)", 0.3092906967350258),
    @("Oil, Gas & Consumable Fuels(125)", 0.3066021482343901),
    @("Construction & Engineering(21)", 0.2901969734324868),
    @("Thrifts & Mortgage Finance(47)", 0.2873668454502063),
    @("Life Sciences Tools & Services(19)", 0.266813741250306),
    @("Insurance(75)", 0.2493363790682562),
    @("Electrical Equipment(28)", 0.2311967513597625),
    @("Specialty Retail(59)", 0.2283518316305553),
    @("Capital Markets(76)", 0.2190033797820798),
    @("ETF(303)", 0.2050584300922346),
    @("Beverages(21)", 0.1984899453511643),
    @("Metals & Mining(106)", 0.188150235078575),
    @("Real Estate Management & Development(23)", 0.182575396346169),
    @("Commercial Services & Supplies(52)", 0.1734925025275525),
    @("Hotels, Restaurants & Leisure(51)", 0.1701264839495237),
    @("Electronic Equipment, Instruments & C...(78)", 0.1698498170234863),
    @("IT Services(52)", 0.1484281979832643),
    @("Health Care Providers & Services(47)", 0.1477258185257507),
    @("Household Durables(39)", 0.1409376841253368),
    @("Software(70)", 0.1389822047452081),
    @("Professional Services(35)", 0.1366943865403495),
    @("Food Products(46)", 0.1263535284096473),
    @("Pharmaceuticals(53)", 0.1196534448946543),
    @("Media(42)", 0.1142760206327241),
    @("Health Care Equipment & Supplies(86)", 0.1030140121109655),
    @("Communications Equipment(45)", 0.09691230739560371),
    @("Equity Real Estate Investment Trusts ...(98)", 0.07417038908594918),
    @("Biotechnology(128)", 0.07117618574307928)
)

for ($i = 0; $i -lt $sectorData.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $sectorData[$i][0]
    $ws.Cells.Item($row, 2).Value = $sectorData[$i][1]
}
